# Update odds values in row 2 of Sheet1 as per the Flashscore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J2").Value = 1.7
$ws.Range("K2").Value = 2.65
$ws.Range("L2").Value = 7.1
$ws.Range("N2").Value = 9.75
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 4.9
$ws.Range("R2").Value = 2.57
$ws.Range("T2").Value = 3.5
$ws.Range("U2").Value = 1.78
$ws.Range("V2").Value = 1.93
$ws.Range("W2").Value = 9.25
$ws.Range("X2").Value = 7.4
$ws.Range("AA2").Value = 10
$ws.Range("AC2").Value = 9.75
$ws.Range("AD2").Value = 11
$ws.Range("AG2").Value = 500
$ws.Range("AM2").Value = 70
$ws.Range("AN2").Value = 3.35
$ws.Range("AO2").Value = 5.6
$ws.Range("AQ2").Value = 13.5
$ws.Range("AR2").Value = 32
$ws.Range("AS2").Value = 150
$ws.Range("AT2").Value = 3.5
$ws.Range("BA2").Value = 300
